$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.648.61'
$ws.Range('E2').Value = '  -2.15%  '
$ws.Range('D3').Value = '1.797.81'
$ws.Range('E3').Value = '  -1.81%  '
$ws.Range('D4').Value = '''1.004'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''231.54'
$ws.Range('D6').Value = '''0.5891'
$ws.Range('E6').Value = '  -2.41%  '
$ws.Range('D7').Value = '''1.003'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').Value = '''0.2762'
$ws.Range('E8').Value = '  -1.11%  '
$ws.Range('D9').Value = '''0.06793'
$ws.Range('E9').Value = '  -3.36%  '
$ws.Range('D10').Value = '''23.19'
$ws.Range('E10').Value = '  -1.30%  '
$ws.Range('D11').Value = '''0.07530'
$ws.Range('E11').Value = '  -1.84%  '
$ws.Range('D12').Value = '1.798.40'
$ws.Range('E12').Value = '  -1.89%  '
$ws.Range('D13').Value = '''4.762'
$ws.Range('D14').Value = '''0.6190'
$ws.Range('E14').Value = '  -1.08%  '
$ws.Range('D15').Value = '2.040.54'
$ws.Range('E15').Value = '  -1.86%  '
$ws.Range('D16').Value = '''0.000009120'
$ws.Range('E16').Value = '  -7.88%  '
$ws.Range('D17').Value = '''75.56'
$ws.Range('E17').Value = '  -4.59%  '
$ws.Range('D18').Value = '28.604.21'
$ws.Range('E18').Value = '  -2.21%  '
$ws.Range('D19').Value = '''5.459'
$ws.Range('D20').Value = '''1.004'
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('D21').Value = '''210.46'
$ws.Range('E21').Value = '  -6.53%  '
$ws.Range('E22').Value = '  -1.79%  '
$ws.Range('D23').Value = '''6.810'
$ws.Range('E23').Value = '  -2.84%  '
$ws.Range('D24').Value = '''1.004'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').Value = '''153.63'
$ws.Range('E25').Value = '  -1.89%  '
$ws.Range('D26').Value = '''7.880'
$ws.Range('E26').Value = '  -1.59%  '
$ws.Range('E27').Value = '  -2.27%  '
$ws.Range('D28').Value = '''16.43'
$ws.Range('E28').Value = '  -0.72%  '
$ws.Range('D29').Value = '''1.424'
$ws.Range('E29').Value = '  -3.74%  '
$ws.Range('D30').Value = '''0.06119'
$ws.Range('E30').Value = '  -1.12%  '
$ws.Range('E31').Value = '  -1.48%  '
$ws.Range('D32').Value = '''3.800'
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('D33').Value = '''3.782'
$ws.Range('E33').Value = '  -1.32%  '
$ws.Range('E34').Value = '  -0.82%  '
$ws.Range('E35').Value = '  -6.01%  '
$ws.Range('D36').Value = '''0.6407'
$ws.Range('E36').Value = '  -0.88%  '
$ws.Range('E37').Value = '  -1.79%  '
$ws.Range('D38').Value = '''2.715'
$ws.Range('E38').Value = '  -0.89%  '
$ws.Range('D39').Value = '''6.538'
$ws.Range('E39').Value = '  -0.19%  '
$ws.Range('D40').Value = '''0.01697'
$ws.Range('E40').Value = '  -2.34%  '
$ws.Range('D41').Value = '1.147.40'
$ws.Range('E41').Value = '  -6.22%  '
$ws.Range('D42').Value = '''0.8857'
$ws.Range('E42').Value = '  -1.44%  '
$ws.Range('D43').Value = '''1.008'
$ws.Range('E43').Value = '  +0.32%  '
$ws.Range('D44').Value = '''100.05'
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('D45').Value = '1.943.55'
$ws.Range('E45').Value = '  -2.25%  '
$ws.Range('D46').Value = '''60.17'
$ws.Range('E46').Value = '  -3.82%  '
$ws.Range('E47').Value = '  -3.59%  '
$ws.Range('D48').Value = '''1.587'
$ws.Range('E48').Value = '  +0.41%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '''0.05454'
$ws.Range('E49').Value = '  -0.98%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''8.290'
$ws.Range('E50').Value = '  -2.89%  '
$ws.Range('D51').Value = '''0.4481'
$ws.Range('E51').Value = '  -1.81%  '
